# Update container names in mactaquac parsers
# - Fix the misaligned helper-text row (row 3) on both the "Ponds" and
#   "Eggrooms" sheets so each comment lines up under its real header.
# - Give the tank/trough comments clearer examples ("E.g. LP1" / "E.g. TR1").
# - Make "Eggrooms" the active/selected tab (it was "Ponds" before).

$wb = $excel.ActiveWorkbook

$wsPonds = $wb.Worksheets.Item("Ponds")
$wsEggrooms = $wb.Worksheets.Item("Eggrooms")

# --- Ponds sheet: fix row 3 helper text so it matches the row 2 headers ---
$wsPonds.Range("D3").Value = "Name of tank. E.g. LP1"
$wsPonds.Range("E3").Value = "Treatment name, must match treatment code in database. Eg. Formaldehyde"
$wsPonds.Range("F3").Value = "Eg. 6"
$wsPonds.Range("G3").Value = "Units can be set to (Gal), (ml) or (kg)."
$wsPonds.Range("J3").Value = "Eg. 1:500"
$wsPonds.Range("K3").Value = "Eg. 1.25"
$wsPonds.Range("L3").Value = "Eg. AB, CD"

# --- Eggrooms sheet: fix row 3 helper text so it matches the row 2 headers ---
$wsEggrooms.Range("D3").Value = "Name of trough. E.g. TR1"
$wsEggrooms.Range("E3").Value = "Treatment name, must match treatment code in database. Eg. Formaldehyde"
$wsEggrooms.Range("F3").Value = "Eg. 6"
$wsEggrooms.Range("G3").Value = "Units can be set to (Gal), (ml) or (kg)."
$wsEggrooms.Range("J3").Value = "Eg. 1:500"
$wsEggrooms.Range("K3").Value = "Eg. 1.25"
$wsEggrooms.Range("L3").Value = "Eg. AB, CD"

# --- Selection / active-sheet bookkeeping to match the saved workbook state ---
$wsPonds.Activate()
$wsPonds.Range("D4").Select()

$wsEggrooms.Activate()
$wsEggrooms.Range("D3").Select()
